# "Add files via upload" commit: the workbook's sheet was renamed and a new
# instructional note was added above the "Tabela da Mahle" block; everything
# else in the upstream diff (fileVersion/rupBuild, absPath, revisionPtr,
# window geometry, calcFeatures ext, dxf/cellXf reordering, theme display
# name, recalculated default column widths, etc.) is Excel-version save
# noise that isn't a deliberate content edit, so we only reproduce the
# substantive changes below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Planilha1" -> "AV - DIAGNÓSTICA"
$ws.Name = "AV - DIAGNÓSTICA"

# New note cell just under the title (becomes shared string index 23)
$ws.Range("B3").Value = "Add 30% no preço de venda."

# Author's last selection before saving moved to E20
$ws.Range("E20").Select() | Out-Null
